# Add a new "2023" column (N) to the CITES permits table, mirroring the
# formatting of the existing "2022" column (M).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 is a thin divider row with no values, just the bottom-border style.
$ws.Range("M3").Copy()
$ws.Range("N3").PasteSpecial(-4122)  # xlPasteFormats

# Row 4 holds the year headers (2013 ... 2022); extend it with 2023.
$ws.Range("M4").Copy()
$ws.Range("N4").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("N4").Value = 2023

# Row 5 holds the data values for each year; add the new 2023 figure.
$ws.Range("M5").Copy()
$ws.Range("N5").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("N5").Value = 553

# Row 3's height becomes an explicit custom height once the new cell is
# added to it (matches Excel's re-flow of the divider row).
$ws.Rows.Item(3).RowHeight = 13.5

# Reset the selection back to the top of the sheet (it had drifted to O4,
# just past the old last column).
$ws.Range("A1").Select()

$excel.CutCopyMode = $false
